$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column BF holds a "Date" column that was recorded using the wrong format
# (e.g. "6-19-2012-13" instead of "2013-06-19"). Fix the 30 data rows
# (rows 2 through 31) to use the correct ISO-like date string.
for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Range("BF$row")
    $cell.NumberFormat = "@"
    $cell.Value = "2013-06-19"
    $cell.Style = "Normal"
}
